# Horas.xlsx - "Timing need to be fixed"
#
# The sheet tracks weekly time entries per person. This edit:
#  - Consolidates the roster: "Martin"/"Miguel"/"Miguel Lolo"/"Raul" rows are
#    dropped, "Luisito cara pito" moves up to row 3 and "Victor" moves up to
#    row 4 (with a corrected time, 0:00 -> 12:30).
#  - Gives "Laia Gonzalez" (row 2) a bunch of additional logged time slots
#    (columns up through AA) and bumps her completed-count (G2) from 19 to 26.
#  - Fixes the SUMA() ranges that total the first and second tables so they
#    cover the right rows after the roster change.
#  - Adds a new coordination entry "asdfasd" (row 65) with the standard 36.28
#    rate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Laia González): extend with more logged time slots, bump G2 ---
$ws.Range("G2").Value = 26
$ws.Range("U2").Value = "2023-01-06 --> 4:00"
$ws.Range("V2").Value = "2023-01-06 --> 0:00"
$ws.Range("W2").Value = "2023-01-12 --> 3:00"
$ws.Range("X2").Value = "2023-01-12 --> 3:00"
$ws.Range("Y2").Value = "2023-01-12 --> 3:00"
$ws.Range("Z2").Value = "2023-01-12 --> 12:00"
$ws.Range("AA2").Value = "2023-01-12 --> 12:00"

# --- Row 3: was "Martín" (now removed), becomes "Luisito cara pito" ---
$ws.Range("A3").Value = "Luisito cara pito"
$ws.Range("B3").Value = "2023-01-12 --> 14:00"
$ws.Range("C3").Value = "2023-01-12 --> 14:00"
$ws.Range("D3").Value = "2023-01-13 --> 15:00"
$ws.Range("E3").Value = "2023-01-13 --> 15:00"
$ws.Range("G3").Value = 4

# --- Row 4: was "Luisito cara pito", becomes "Victor" with fixed time ---
$ws.Range("A4").Value = "Victor"
$ws.Range("B4").Value = "2023-01-12 --> 12:30"
$ws.Range("C4:E4").ClearContents()
$ws.Range("G4").Value = 1

# --- Rows 5-8 ("Miguel", "Miguel Lolo", "Raul", old "Victor"): removed ---
$ws.Range("A5:G8").ClearContents()

# --- Fix the SUMA() ranges now that the roster rows shifted ---
$ws.Range("G38").Formula = "=SUMA(G2:G37)"
$ws.Range("G53").Formula = "=SUMA(G42:G52)"

# --- New coordination entry row ---
$ws.Range("A65").Value = "asdfasd"
$ws.Range("B65").Value = 36.28
